# Actualización automática 2025-07-29 13:30:09
# Applies a new 24.39 PORCELANATO sale for DANIELA ELIZABETH BECERRA BECERRA
# (OFICINA-CATAECSA) recorded in July, plus the resulting totals/rollups.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": per-client sales by product group ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M7").Value = 24.39
$wsGrupo.Range("M19").Value = "2 de 17"

# --- Sheet "VENTA MENSUAL": per-client sales by month ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F7").Value = 24.39
$wsMensual.Range("F19").Value = 1644.26

# --- Sheet "CUMPLIMIENTO MENSUAL": compliance summary totals ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PORCELANATO row
$wsCumplimiento.Range("D3").Value = 14220.05
$wsCumplimiento.Range("E3").Value = -496.7099999999991
$wsCumplimiento.Range("F3").Value = 1.036194541562039

# TOTAL row
$wsCumplimiento.Range("D4").Value = 20593.26
$wsCumplimiento.Range("E4").Value = -6869.919999999999
$wsCumplimiento.Range("F4").Value = 1.500601165605457
